# Generate Report for Handback
#
# The b67d0044-7a31-4bb6-ae28-fafd7098fe29.md file has finished its
# handback cycle: mark it as "Handed back: in sync with en-US" (it was
# "Ready for handoff") on the Overview sheet and on each language sheet,
# and stamp the per-language "Latest Handback DateTime" with the handback
# timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for b67d0044...md is row 3 ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for b67d0044...md is row 3 ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-03-10 16:44:25"

# --- de-de sheet: row for b67d0044...md is row 3 ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2016-03-10 16:44:34"
